$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.555.58'
$ws.Range('E2').Value = '  -3.72%  '
$ws.Range('D3').Value = '2.510.38'
$ws.Range('E3').Value = '  -5.14%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''577.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('D6').Value = '''167.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.67%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').Value = '2.509.53'
$ws.Range('E9').Value = '  -5.12%  '
$ws.Range('E10').Value = '  -7.07%  '
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('E12').Value = '  -4.14%  '
$ws.Range('D13').Value = '''4.87'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '2.964.57'
$ws.Range('E14').Value = '  -5.31%  '
$ws.Range('D15').Value = '69.498.45'
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('D16').Value = '''0.0000175'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.96%  '
$ws.Range('D17').Value = '''24.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.06%  '
$ws.Range('D18').Value = '2.514.77'
$ws.Range('E18').Value = '  -4.39%  '
$ws.Range('D19').Value = '''11.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.57%  '
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('D21').Value = '''350.68'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('E22').Value = '  -4.88%  '
$ws.Range('D23').Value = '''1.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.05%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '''69.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.24%  '
$ws.Range('D26').Value = '''4.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.00%  '
$ws.Range('D27').Value = '''9.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.02%  '
$ws.Range('D28').Value = '2.642.46'
$ws.Range('D29').Value = '''0.996'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('D30').Value = '0.0₃0904'
$ws.Range('E30').Value = '  -5.40%  '
$ws.Range('D31').Value = '''7.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('D32').Value = '''478.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('D33').Value = '''1.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = '''1.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '''154.12'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.88%  '
$ws.Range('D39').Value = '''18.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.13%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '''4.76'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.18%  '
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('D43').Value = '''1.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.33%  '
$ws.Range('E44').Value = '  -13.57%  '
$ws.Range('E45').Value = '  -8.71%  '
$ws.Range('D46').Value = '''38.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').Value = '''144.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.31%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').Value = '''3.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '''0.531'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.69%  '
$ws.Range('D50').Value = '''1.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.25%  '
$ws.Range('E51').Value = '  -2.41%  '
